# Update "Avvisi" (notices) management REST endpoints in the
# "completezza funzionale" tracking sheet.
#
# Endpoints added:
#   /avviso/crea
#   /avviso/segna-come-letto/{id_avviso}
#   /avviso/segna-come-nascosto/{id_avviso}
#
# The previous placeholder row "setAvvisoViewed #incompleto" is replaced
# by real endpoints, the existing "Avvisi Management Requests" endpoint
# list (column G) is shifted down one row to make room for
# "/avvisi?id_ristorante=" at the top, and a new (still empty) row is
# started below, marked with an underlined font as a placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 / D4 previously carried a cell style (s="13") that applied the
# default font with no visible effect - restore the plain/default look.
$ws.Range("B3").ClearFormats()
$ws.Range("D4").ClearFormats()

# Rebuild the "Avvisi Management Requests" endpoint list (column G,
# rows 3-8) with the up-to-date set of REST routes.
$ws.Range("G3").Value = "/avvisi?id_ristorante="
$ws.Range("G4").Value = "/avvisi-hidden/{id_user}"
$ws.Range("G5").Value = "/avvisi-viewed/{id_user}"
$ws.Range("G6").Value = "/avviso/segna-come-letto/{id_avviso}"
$ws.Range("G7").Value = "/avviso/segna-come-nascosto/{id_avviso}"
$ws.Range("G8").Value = "/avviso/crea"

# New placeholder entry started at G10 - underlined, still empty.
$ws.Range("G10").Font.Underline = $true

# Scroll the sheet one column to the right and move the selection to
# the newly added placeholder cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("G10").Select()
